$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 50, shifting existing rows 50-74 down to 53-77
$ws.Rows.Item(50).Resize(3).Insert()

# Copy the date cell's number format from what is now row 53 (original row 50) into the new rows
$ws.Range("D50:D52").NumberFormat = $ws.Range("D53").NumberFormat

$newDate = [DateTime]::FromOADate(44488)

# Row 50 - Especial
$ws.Cells.Item(50, 1).Value = 8
$ws.Cells.Item(50, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value = "Coquimbo"
$ws.Cells.Item(50, 4).Value = $newDate
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100107
$ws.Cells.Item(50, 8).Value = "Otros"
$ws.Cells.Item(50, 9).Value = 100107002
$ws.Cells.Item(50, 10).Value = "Chirimoya"
$ws.Cells.Item(50, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(50, 12).Value = "Especial"
$ws.Cells.Item(50, 13).Value = 300
$ws.Cells.Item(50, 14).Value = 2200
$ws.Cells.Item(50, 15).Value = 2300
$ws.Cells.Item(50, 16).Value = 2250
$ws.Cells.Item(50, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(50, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(50, 19).Value = 2250
$ws.Cells.Item(50, 20).Value = 1

# Row 51 - Primera
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = $newDate
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100107
$ws.Cells.Item(51, 8).Value = "Otros"
$ws.Cells.Item(51, 9).Value = 100107002
$ws.Cells.Item(51, 10).Value = "Chirimoya"
$ws.Cells.Item(51, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 300
$ws.Cells.Item(51, 14).Value = 1900
$ws.Cells.Item(51, 15).Value = 2000
$ws.Cells.Item(51, 16).Value = 1950
$ws.Cells.Item(51, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(51, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 19).Value = 1950
$ws.Cells.Item(51, 20).Value = 1

# Row 52 - Segunda
$ws.Cells.Item(52, 1).Value = 8
$ws.Cells.Item(52, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = $newDate
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100107
$ws.Cells.Item(52, 8).Value = "Otros"
$ws.Cells.Item(52, 9).Value = 100107002
$ws.Cells.Item(52, 10).Value = "Chirimoya"
$ws.Cells.Item(52, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 240
$ws.Cells.Item(52, 14).Value = 1400
$ws.Cells.Item(52, 15).Value = 1500
$ws.Cells.Item(52, 16).Value = 1450
$ws.Cells.Item(52, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(52, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(52, 19).Value = 1450
$ws.Cells.Item(52, 20).Value = 1
